# Updated cryptos list (Price / Volume(1h) refresh, plus a Kaspa/Mantle row
# swap) to match the latest GitHub Actions scrape.
# Note: several "Price" values are plain decimal numbers (e.g. "601.07").
# Excel would auto-convert those to numeric cells (dropping formatting such
# as trailing zeros) unless we force text entry via a leading apostrophe,
# since the source data stores prices as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.073.63'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '3.736.24'
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''601.07'
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("D6").Value = '''167.21'
$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("D7").Value = '3.734.53'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +0.98%  '

$ws.Range("E10").Value = '  +3.81%  '

$ws.Range("E11").Value = '  +0.39%  '

$ws.Range("D12").Value = '''0.460'
$ws.Range("E12").Value = '  +0.38%  '

$ws.Range("D13").Value = '''38.06'
$ws.Range("E13").Value = '  +0.11%  '

$ws.Range("D14").Value = '''0.0000247'
$ws.Range("E14").Value = '  +1.57%  '

$ws.Range("D15").Value = '4.365.61'
$ws.Range("E15").Value = '  +0.08%  '

$ws.Range("D16").Value = '3.741.84'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").Value = '69.008.19'
$ws.Range("E17").Value = '  +0.37%  '

$ws.Range("D18").Value = '''7.35'
$ws.Range("E18").Value = '  +1.42%  '

$ws.Range("D19").Value = '''17.33'
$ws.Range("E19").Value = '  +0.34%  '

$ws.Range("D20").Value = '''0.113'
$ws.Range("E20").Value = '  -1.67%  '

$ws.Range("D21").Value = '''11.23'
$ws.Range("E21").Value = '  +11.10%  '

$ws.Range("D22").Value = '''491.51'
$ws.Range("E22").Value = '  -1.08%  '

$ws.Range("D23").Value = '''0.726'
$ws.Range("E23").Value = '  +0.43%  '

$ws.Range("E24").Value = '  +7.73%  '

$ws.Range("D25").Value = '''84.75'
$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("D26").Value = '''2.28'
$ws.Range("E26").Value = '  -0.79%  '

$ws.Range("D27").Value = '''12.27'
$ws.Range("E27").Value = '  -0.62%  '

$ws.Range("E28").Value = '  -0.18%  '

$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("E30").Value = '  +1.15%  '

$ws.Range("D31").Value = '''8.21'
$ws.Range("E31").Value = '  +3.56%  '

$ws.Range("D33").Value = '''31.43'
$ws.Range("E33").Value = '  -0.76%  '

$ws.Range("D34").Value = '3.883.11'
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("D36").Value = '3.669.84'
$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '''0.140'
$ws.Range("E38").Value = '  +5.77%  '

$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").Value = '''1.01'
$ws.Range("E39").Value = '  -0.29%  '

$ws.Range("D40").Value = '''5.93'
$ws.Range("E40").Value = '  +1.94%  '

$ws.Range("D41").Value = '''0.325'
$ws.Range("E41").Value = '  +0.20%  '

$ws.Range("E42").Value = '  +5.72%  '

$ws.Range("D43").Value = '''48.78'
$ws.Range("E43").Value = '  -0.48%  '

$ws.Range("D44").Value = '''1.98'
$ws.Range("E44").Value = '  +0.59%  '

$ws.Range("D45").Value = '''423.28'
$ws.Range("E45").Value = '  -2.42%  '

$ws.Range("D46").Value = '''8.46'
$ws.Range("E46").Value = '  +0.53%  '

$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("D48").Value = '''39.98'
$ws.Range("E48").Value = '  -1.37%  '

$ws.Range("D49").Value = '''141.40'
$ws.Range("E49").Value = '  -0.17%  '

$ws.Range("D50").Value = '2.778.49'
$ws.Range("E50").Value = '  +1.27%  '

$ws.Range("D51").Value = '''0.0353'
$ws.Range("E51").Value = '  +0.08%  '
